$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Sprint 6 stand-up update: refresh the "This Week" block (rows 14-16,
# columns B/D) with the latest per-person status text.
$ws.Range("B14").Value = "1) Did changes on ER diagram"
$ws.Range("C14").Value = "1)Prepared Workshop resentation slides on Firebase database and learned queries and architecture."
$ws.Range("D14").Value = "1) Finalaized ER diagram and learned about Android studio and how to integrate Firebase Database in Android studio."

$ws.Range("B15").Value = "2) I will prepare workshop presentation document on Firebase databse."
$ws.Range("C15").Value = "2) I will give presentaion on Firebase database and will do changes in ER diagram."
$ws.Range("D15").Value = "2) I will change the Database based on project requirement."

$ws.Range("B16").Value = "3) NA"
$ws.Range("C16").Value = "3) NA"
$ws.Range("D16").Value = "3) We are finalizing the database."

# Row heights grew to fit the longer wrapped text.
$ws.Rows.Item(14).RowHeight = 62
$ws.Rows.Item(15).RowHeight = 46.5

# Leave the sheet selected where the author last left it.
[void]$ws.Range("C19").Select()
